# Roll the 12-month PL window forward from (Jan 2023 - Dec 2023) to
# (Oct 2023 - Sep 2024): the last three actual months (Oct/Nov/Dec 2023)
# slide into the first three columns, everything else resets to 0 since
# those future months haven't been synced from the AI chat yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: B1:M1 month labels -------------------------------------
# Force text so month/year labels like "Oct 2023" aren't auto-parsed into
# date serials by Excel's input-type inference (mirrors setting the cell
# format to Text before typing the value).
$headers = @("Oct 2023","Nov 2023","Dec 2023","Jan 2024","Feb 2024","Mar 2024","Apr 2024","May 2024","Jun 2024","Jul 2024","Aug 2024","Sep 2024")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
}

# --- Row 2: Cost of Sales - Cost of sales --------------------------------
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 3250
$ws.Range("D2").Value = 1125
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0

# --- Row 3: Gross Profit --------------------------------------------------
$ws.Range("B3").Value = 3213.75
$ws.Range("C3").Value = 1412.5
$ws.Range("D3").Value = 10542.36
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0

# --- Row 4: Expenses - Depreciation Expense -------------------------------
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 99.98999999999999
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0

# --- Row 5: Net Income -----------------------------------------------------
$ws.Range("B5").Value = 1589.76
$ws.Range("C5").Value = -537.8200000000001
$ws.Range("D5").Value = 8068.03
$ws.Range("E5").Value = -1622.41
$ws.Range("F5").Value = -1349.99
$ws.Range("G5").Value = -3015
$ws.Range("H5").Value = -95
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
